$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price update: a new week's Perejil (parsley) record is published at the
# top of the data block, so every existing row from 52 downward shifts down by
# one row and the freshly-reported week's figures land in the now-empty row 52.
$ws.Rows("52").Insert()

$ws.Range("A52").Value = 8
$ws.Range("B52").Value = "Terminal La Palmera de La Serena"
$ws.Range("C52").Value = "Coquimbo"
$ws.Range("D52").Value = 44455
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 100112044
$ws.Range("G52").Value = "Perejil"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 3200
$ws.Range("K52").Value = 2000
$ws.Range("L52").Value = 2500
$ws.Range("M52").Value = 2250
$ws.Range("N52").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O52").Value = "Provincia del Elquí"
$ws.Range("P52").Value = 1500
$ws.Range("Q52").Value = 1.5
$ws.Range("R52").Value = "Hortaliza"
